# Update the "dSF" column (F) values for select rows after repulling data.
# These rows correspond to data index A = 0, 2, 3, 6, 9 (sheet rows 2, 4, 5, 8, 11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -6
$ws.Range("F4").Value = -11
$ws.Range("F5").Value = -3
$ws.Range("F8").Value = -10
$ws.Range("F11").Value = -5
